$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 24
$ws.Range("C2").Value = 9
$ws.Range("E2").Value = " Liam Livingstone"
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = "Caught"
$ws.Range("B3").Value = 13
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = "Bowled"
$ws.Range("E3").Value = " Adil Rashid"
$ws.Range("K3").Value = 32
$ws.Range("L3").Value = 11
$ws.Range("M3").Value = "Caught"
$ws.Range("N3").Value = " Shadab Khan"
$ws.Range("B4").Value = 7
$ws.Range("C4").Value = 4
$ws.Range("E4").Value = " Chris Jordan"
$ws.Range("K4").Value = 32
$ws.Range("L4").Value = 12
$ws.Range("B5").Value = 29
$ws.Range("C5").Value = 7
$ws.Range("D5").Value = "LBW"
$ws.Range("E5").Value = " Chris Woakes"
$ws.Range("K5").Value = 5
$ws.Range("L5").Value = 2
$ws.Range("N5").Value = " Imad Wasim"
$ws.Range("B6").Value = 17
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = "LBW"
$ws.Range("E6").Value = " Liam Livingstone"
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 1
$ws.Range("N6").Value = " Imad Wasim"
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = "LBW"
$ws.Range("E7").Value = " Liam Livingstone"
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = "NOT OUT"
$ws.Range("N7").Value = " "
$ws.Range("B8").Value = 13
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = "LBW"
$ws.Range("E8").Value = " Mark Wood"
$ws.Range("K8").Value = 5
$ws.Range("L8").Value = 2
$ws.Range("N8").Value = " Haris Rauf"
$ws.Range("B9").Value = 21
$ws.Range("C9").Value = 7
$ws.Range("D9").Value = "Caught"
$ws.Range("E9").Value = " Chris Jordan"
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 3
$ws.Range("N9").Value = " Haris Rauf"
$ws.Range("B10").Value = 9
$ws.Range("D10").Value = "Bowled"
$ws.Range("E10").Value = " Chris Woakes"
$ws.Range("K10").Value = 14
$ws.Range("M10").Value = "LBW"
$ws.Range("N10").Value = " Shaheen Afridi"
$ws.Range("B11").Value = 10
$ws.Range("C11").Value = 4
$ws.Range("D11").Value = "LBW"
$ws.Range("E11").Value = " Liam Livingstone"
$ws.Range("K11").Value = 5
$ws.Range("L11").Value = 2
$ws.Range("M11").Value = "Caught"
$ws.Range("N11").Value = " Shaheen Afridi"
$ws.Range("B12").Value = 3
$ws.Range("K12").Value = 4
$ws.Range("L12").Value = 3
$ws.Range("M12").Value = "Caught"
$ws.Range("N12").Value = " Hasan Ali"
$ws.Range("A16").Value = 152
$ws.Range("C16").Value = "'9.4"
$ws.Range("D16").Value = 58
$ws.Range("J16").Value = 104
$ws.Range("L16").Value = "'7.3"
$ws.Range("M16").Value = 45
$ws.Range("A21").Value = "Mark Wood"
$ws.Range("C21").Value = 39
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 19.5
$ws.Range("J21").Value = "Imad Wasim"
$ws.Range("K21").Value = "'1.0"
$ws.Range("L21").Value = 11
$ws.Range("M21").Value = 3
$ws.Range("N21").Value = 11
$ws.Range("A22").Value = "Adil Rashid"
$ws.Range("C22").Value = 27
$ws.Range("E22").Value = 13.5
$ws.Range("J22").Value = "Shadab Khan"
$ws.Range("K22").Value = "'1.0"
$ws.Range("L22").Value = 23
$ws.Range("M22").Value = 1
$ws.Range("N22").Value = 23
$ws.Range("A23").Value = "Chris Jordan"
$ws.Range("C23").Value = 31
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 15.5
$ws.Range("J23").Value = "Haris Rauf"
$ws.Range("L23").Value = 19
$ws.Range("M23").Value = 3
$ws.Range("N23").Value = 9.5
$ws.Range("A24").Value = "Chris Woakes"
$ws.Range("B24").Value = "'2.0"
$ws.Range("C24").Value = 38
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = 19
$ws.Range("J24").Value = "Shaheen Afridi"
$ws.Range("L24").Value = 31
$ws.Range("N24").Value = 15.5
$ws.Range("A25").Value = "Liam Livingstone"
$ws.Range("B25").Value = "'1.4"
$ws.Range("C25").Value = 17
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 12.14
$ws.Range("J25").Value = "Hasan Ali"
$ws.Range("K25").Value = "'1.3"
$ws.Range("L25").Value = 20
$ws.Range("N25").Value = 15.38

$ws.Range("C16").Style = "Normal"
$ws.Range("L16").Style = "Normal"
$ws.Range("K21").Style = "Normal"
$ws.Range("K22").Style = "Normal"
$ws.Range("B24").Style = "Normal"
$ws.Range("B25").Style = "Normal"
$ws.Range("K25").Style = "Normal"
